# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-03-06 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-03-07 Friday", 2)

# Update the answer grid. Only the 5 populated rows (1, 5, 9, 13, 17) of the
# 20x5 table carry text; the others are blank spacer rows.
$t = $d.Tables.Item(1)

$answers = @{
    1  = @("24÷5=4, 4", "41÷2=20, 1", "75÷9=8, 3", "68÷8=8, 4", "32÷8=4, 0")
    5  = @("21÷9=2, 3", "33÷9=3, 6", "77÷5=15, 2", "45÷8=5, 5", "72÷5=14, 2")
    9  = @("67÷9=7, 4", "68÷7=9, 5", "49÷7=7, 0", "50÷5=10, 0", "10÷6=1, 4")
    13 = @("62÷4=15, 2", "20÷3=6, 2", "83÷3=27, 2", "67÷7=9, 4", "97÷3=32, 1")
    17 = @("62÷8=7, 6", "69÷8=8, 5", "56÷3=18, 2", "85÷7=12, 1", "82÷9=9, 1")
}

foreach ($row in $answers.Keys) {
    $values = $answers[$row]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
